# Update the "Förändrad" (Changed) date column (C) from 45233 to 45243
# for all data rows (2 through 15) on the active worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 15; $row++) {
    $ws.Cells.Item($row, 3).Value = 45243
}
